$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 3 (PORCELANATO): update VENTA, POR CUMPLIR, CUMPLIMIENTO
$ws.Range("D3").Value = 1166.29
$ws.Range("E3").Value = 12557.05
$ws.Range("F3").Value = 0.08498587078655778

# Row 4 (TOTAL): update VENTA, POR CUMPLIR, CUMPLIMIENTO
$ws.Range("D4").Value = 4758.96
$ws.Range("E4").Value = 8964.379999999999
$ws.Range("F4").Value = 0.3467785539088881
